$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Iteration_SheetOne")

# Add hyperlinks for the dropdown-value keyword columns (username/email and password)
$null = $ws.Hyperlinks.Add($ws.Range("G2"), "mailto:karthik.sharma041992@gmail.com")
$null = $ws.Hyperlinks.Add($ws.Range("G3"), "mailto:karthik.sharma041992@gmail.com")
$null = $ws.Hyperlinks.Add($ws.Range("H2"), "mailto:karthik.sharma041992@gmail.com")
$null = $ws.Hyperlinks.Add($ws.Range("H3"), "mailto:karthik.sharma041992@gmail.com")

# Set the displayed values
$ws.Range("G2").Value = "karthik.sharma041992@gmail.com"
$ws.Range("G3").Value = "karthik.sharma041992@gmail.com"
$ws.Range("H2").Value = "April@2018"
$ws.Range("H3").Value = "April@2018"

# Keep the hyperlink look consistent with the existing hyperlink cells
$ws.Range("G2:H3").Style = "Hyperlink"

# Update the active selection to match the final cursor position
$null = $ws.Range("H3").Select()
